# Auto-generated edit script: update country COVID stats and reorder/rename rows per source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 29 de Marzo de 2020 a las 18:50"

# Row 4
$ws.Cells.Item(4,2).Value = 131403
$ws.Cells.Item(4,3).Value = 7825
$ws.Cells.Item(4,5).Value = 124696
$ws.Cells.Item(4,7).Value = 108
$ws.Cells.Item(4,8).Value = 2329

# Row 8
$ws.Cells.Item(8,2).Value = 60659
$ws.Cells.Item(8,3).Value = 2964
$ws.Cells.Item(8,4).Value = 9211
$ws.Cells.Item(8,5).Value = 50966
$ws.Cells.Item(8,6).Value = 1979
$ws.Cells.Item(8,7).Value = 49
$ws.Cells.Item(8,8).Value = 482

# Row 16
$ws.Cells.Item(16,1).Value = "Turquia"
$ws.Cells.Item(16,2).Value = 9217
$ws.Cells.Item(16,3).Value = 1815
$ws.Cells.Item(16,4).Value = 105
$ws.Cells.Item(16,5).Value = 8981
$ws.Cells.Item(16,6).Value = 568
$ws.Cells.Item(16,7).Value = 23
$ws.Cells.Item(16,8).Value = 131

# Row 17
$ws.Cells.Item(17,1).Value = "Austria"
$ws.Cells.Item(17,2).Value = 8672
$ws.Cells.Item(17,3).Value = 401
$ws.Cells.Item(17,4).Value = 479
$ws.Cells.Item(17,5).Value = 8107
$ws.Cells.Item(17,6).Value = 187
$ws.Cells.Item(17,7).Value = 18
$ws.Cells.Item(17,8).Value = 86

# Row 19
$ws.Cells.Item(19,2).Value = 5886
$ws.Cells.Item(19,3).Value = 231
$ws.Cells.Item(19,5).Value = 5315

# Row 21
$ws.Cells.Item(21,2).Value = 3980
$ws.Cells.Item(21,3).Value = 345
$ws.Cells.Item(21,5).Value = 3738

# Row 25
$ws.Cells.Item(25,2).Value = 2743
$ws.Cells.Item(25,3).Value = 112
$ws.Cells.Item(25,5).Value = 2719

# Row 31
$ws.Cells.Item(31,2).Value = 1890
$ws.Cells.Item(31,3).Value = 67
$ws.Cells.Item(31,5).Value = 1830
$ws.Cells.Item(31,7).Value = 9
$ws.Cells.Item(31,8).Value = 57

# Row 33
$ws.Cells.Item(33,5).Value = 1549
$ws.Cells.Item(33,7).Value = 5
$ws.Cells.Item(33,8).Value = 42

# Row 61
$ws.Cells.Item(61,1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(61,2).Value = 570
$ws.Cells.Item(61,3).Value = 102
$ws.Cells.Item(61,4).Value = 58
$ws.Cells.Item(61,5).Value = 510
$ws.Cells.Item(61,6).Value = 2
$ws.Cells.Item(61,8).Value = 2

# Row 62
$ws.Cells.Item(62,1).Value = "Irak"
$ws.Cells.Item(62,2).Value = 547
$ws.Cells.Item(62,3).Value = 41
$ws.Cells.Item(62,4).Value = 143
$ws.Cells.Item(62,5).Value = 362
$ws.Cells.Item(62,6).Value = 0
$ws.Cells.Item(62,8).Value = 42

# Row 63
$ws.Cells.Item(63,1).Value = "Nueva Zelanda"
$ws.Cells.Item(63,2).Value = 514
$ws.Cells.Item(63,3).Value = 0
$ws.Cells.Item(63,4).Value = 56
$ws.Cells.Item(63,5).Value = 457
$ws.Cells.Item(63,6).Value = 1
$ws.Cells.Item(63,8).Value = 1

# Row 64
$ws.Cells.Item(64,1).Value = "Argelia"
$ws.Cells.Item(64,2).Value = 511
$ws.Cells.Item(64,3).Value = 57
$ws.Cells.Item(64,4).Value = 31
$ws.Cells.Item(64,5).Value = 449
$ws.Cells.Item(64,6).Value = 0
$ws.Cells.Item(64,7).Value = 2
$ws.Cells.Item(64,8).Value = 31

# Row 65
$ws.Cells.Item(65,1).Value = "Barein"
$ws.Cells.Item(65,2).Value = 499
$ws.Cells.Item(65,3).Value = 23
$ws.Cells.Item(65,4).Value = 272
$ws.Cells.Item(65,5).Value = 223
$ws.Cells.Item(65,6).Value = 1
$ws.Cells.Item(65,8).Value = 4

# Row 91
$ws.Cells.Item(91,4).Value = 25
$ws.Cells.Item(91,5).Value = 163

# Row 97
$ws.Cells.Item(97,1).Value = "Uzbekistan"
$ws.Cells.Item(97,2).Value = 144
$ws.Cells.Item(97,3).Value = 40
$ws.Cells.Item(97,4).Value = 7
$ws.Cells.Item(97,5).Value = 135
$ws.Cells.Item(97,6).Value = 8
$ws.Cells.Item(97,8).Value = 2

# Row 98
$ws.Cells.Item(98,1).Value = "Senegal"
$ws.Cells.Item(98,2).Value = 142
$ws.Cells.Item(98,3).Value = 12
$ws.Cells.Item(98,4).Value = 27
$ws.Cells.Item(98,5).Value = 115

# Row 99
$ws.Cells.Item(99,1).Value = "Costa de Marfil"
$ws.Cells.Item(99,2).Value = 140
$ws.Cells.Item(99,3).Value = 39
$ws.Cells.Item(99,4).Value = 3
$ws.Cells.Item(99,5).Value = 137
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,8).Value = 0

# Row 100
$ws.Cells.Item(100,1).Value = "Cuba"
$ws.Cells.Item(100,2).Value = 139
$ws.Cells.Item(100,3).Value = 20
$ws.Cells.Item(100,4).Value = 4
$ws.Cells.Item(100,5).Value = 132
$ws.Cells.Item(100,6).Value = 2
$ws.Cells.Item(100,8).Value = 3

# Row 101
$ws.Cells.Item(101,1).Value = "Brunei"
$ws.Cells.Item(101,2).Value = 126
$ws.Cells.Item(101,3).Value = 6
$ws.Cells.Item(101,4).Value = 34
$ws.Cells.Item(101,5).Value = 91
$ws.Cells.Item(101,6).Value = 1
$ws.Cells.Item(101,8).Value = 1

# Row 102
$ws.Cells.Item(102,1).Value = "Afganistan"
$ws.Cells.Item(102,2).Value = 120
$ws.Cells.Item(102,3).Value = 10
$ws.Cells.Item(102,4).Value = 2
$ws.Cells.Item(102,5).Value = 114
$ws.Cells.Item(102,6).Value = 0
$ws.Cells.Item(102,8).Value = 4

# Row 146
$ws.Cells.Item(146,1).Value = "Mali"
$ws.Cells.Item(146,3).Value = 0

# Row 147
$ws.Cells.Item(147,1).Value = "Niger"
$ws.Cells.Item(147,3).Value = 8

# Row 161
$ws.Cells.Item(161,1).Value = "Birmania"
$ws.Cells.Item(161,3).Value = 2
$ws.Cells.Item(161,4).Value = 0
$ws.Cells.Item(161,5).Value = 10

# Row 162
$ws.Cells.Item(162,1).Value = "Groenlandia"
$ws.Cells.Item(162,2).Value = 10
$ws.Cells.Item(162,4).Value = 2
$ws.Cells.Item(162,5).Value = 8

# Row 164
$ws.Cells.Item(164,1).Value = "Granada"
$ws.Cells.Item(164,3).Value = 2
$ws.Cells.Item(164,5).Value = 9
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 0

# Row 165
$ws.Cells.Item(165,1).Value = "Siria"
$ws.Cells.Item(165,2).Value = 9
$ws.Cells.Item(165,3).Value = 4
$ws.Cells.Item(165,7).Value = 1
$ws.Cells.Item(165,8).Value = 1

# Row 166
$ws.Cells.Item(166,1).Value = "Laos"

# Row 167
$ws.Cells.Item(167,1).Value = "Seychelles"
